$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.372.10'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.81%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.685.90'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.35%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '683.74'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.44'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -6.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.685.28'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.32%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -5.96%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -8.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.20'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.437'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -10.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000233'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -6.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.306.95'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '32.57'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -11.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.683.42'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.408.59'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '15.88'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -9.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.45'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -10.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '474.02'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -7.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.89'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.649'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -9.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.64'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.829.73'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.22%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -11.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.93'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -13.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.19'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -10.91%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.69'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -10.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.75'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -13.28%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -10.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.68'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -9.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.997'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '26.71'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -8.48%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -7.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.21'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -11.98%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -7.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.27'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.57%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0906'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -10.38%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.941'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -7.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '165.62'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.98'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.73'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -15.23%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.78'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.97%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.31'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -6.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.11'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.53%  '
$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000275'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -9.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.88'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -9.21%  '
